$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/7/2025  Through  7/13/2025"

# --- Crime Complaints table updates (rows 14-31) ---
# Cells whose style/category is unchanged: just update the value
$ws.Range("N14").Value = -94.736842105263
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 9
$ws.Range("J15").Value = 14
$ws.Range("K15").Value = -35.714285714285
$ws.Range("L15").Value = -10
$ws.Range("M15").Value = 28.571428571428
$ws.Range("C16").Value = 3
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 37.5
$ws.Range("I16").Value = 64
$ws.Range("J16").Value = 45
$ws.Range("K16").Value = 42.222222222222
$ws.Range("L16").Value = -4.477611940298
$ws.Range("M16").Value = -55.862068965517
$ws.Range("N16").Value = -92.784667418263
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 300
$ws.Range("F17").Value = 34
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = 100
$ws.Range("I17").Value = 201
$ws.Range("J17").Value = 139
$ws.Range("K17").Value = 44.604316546762
$ws.Range("L17").Value = 14.857142857142
$ws.Range("M17").Value = 21.818181818181
$ws.Range("N17").Value = -59.146341463414
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = -40
$ws.Range("I18").Value = 36
$ws.Range("J18").Value = 39
$ws.Range("K18").Value = -7.692307692307
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -77.070063694267
$ws.Range("N18").Value = -96.449704142011
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 300
$ws.Range("F19").Value = 21
$ws.Range("G19").Value = 15
$ws.Range("H19").Value = 40
$ws.Range("I19").Value = 163
$ws.Range("J19").Value = 128
$ws.Range("K19").Value = 27.34375
$ws.Range("L19").Value = -27.87610619469
$ws.Range("M19").Value = -21.634615384615
$ws.Range("N19").Value = -51.632047477744
$ws.Range("C20").Value = 4
$ws.Range("E20").Value = 300
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 83.333333333333
$ws.Range("I20").Value = 51
$ws.Range("J20").Value = 42
$ws.Range("K20").Value = 21.428571428571
$ws.Range("L20").Value = -7.272727272727
$ws.Range("M20").Value = -39.285714285714
$ws.Range("N20").Value = -92.776203966005
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 9
$ws.Range("E21").Value = 166.666666666667
$ws.Range("G21").Value = 53
$ws.Range("H21").Value = 52.830188679245
$ws.Range("I21").Value = 525
$ws.Range("J21").Value = 411
$ws.Range("K21").Value = 27.737226277372
$ws.Range("L21").Value = -8.216783216783
$ws.Range("M21").Value = -31.640625
$ws.Range("N21").Value = -85
$ws.Range("I22").Value = 2
$ws.Range("K22").Value = -71.428571428571
$ws.Range("L22").Value = -81.818181818181
$ws.Range("M22").Value = -92.592592592592
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = -8.333333333333
$ws.Range("F24").Value = 74
$ws.Range("G24").Value = 60
$ws.Range("H24").Value = 23.333333333333
$ws.Range("I24").Value = 520
$ws.Range("J24").Value = 456
$ws.Range("K24").Value = 14.035087719298
$ws.Range("L24").Value = -10.652920962199
$ws.Range("M24").Value = 3.79241516966
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -60
$ws.Range("F25").Value = 9
$ws.Range("G25").Value = 9
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 133
$ws.Range("J25").Value = 83
$ws.Range("K25").Value = 60.240963855421
$ws.Range("L25").Value = 9.9173553719
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = -21.428571428571
$ws.Range("F26").Value = 47
$ws.Range("G26").Value = 53
$ws.Range("H26").Value = -11.320754716981
$ws.Range("I26").Value = 263
$ws.Range("J26").Value = 297
$ws.Range("K26").Value = -11.447811447811
$ws.Range("L26").Value = -1.498127340823
$ws.Range("M26").Value = -40.227272727272
$ws.Range("E27").Value = 0
$ws.Range("H27").Value = -66.666666666666
$ws.Range("I27").Value = 9
$ws.Range("J27").Value = 17
$ws.Range("K27").Value = -47.058823529411
$ws.Range("L27").Value = -25
$ws.Range("C28").Value = 3
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 33
$ws.Range("K28").Value = -15.384615384615
$ws.Range("L28").Value = 13.793103448275
$ws.Range("L29").Value = -57.142857142857
$ws.Range("M29").Value = -83.333333333333
$ws.Range("N29").Value = -95.588235294117
$ws.Range("L30").Value = -50
$ws.Range("M30").Value = -78.571428571428
$ws.Range("N30").Value = -95.522388059701
$ws.Range("L31").Value = 42.857142857142

# Cells that change category (text <-> numeric): copy number format from a
# same-category reference cell (row 33, untouched by this edit) then set the value.
# For new text values, prefix with an apostrophe so the numeric-looking string
# ("0") is stored as text, matching the shared-string placeholders used elsewhere.
$ws.Range("C15").Value = 1
$ws.Range("I33").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D15").Value = 1
$ws.Range("I33").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = 0
$ws.Range("L33").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("F15").Value = 1
$ws.Range("I33").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("C18").Value = "'0"
$ws.Range("C33").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D18").Value = 1
$ws.Range("I33").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").Value = -100
$ws.Range("L33").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("D22").Value = "'0"
$ws.Range("C33").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "'***.*"
$ws.Range("C33").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("G22").Value = "'0"
$ws.Range("C33").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("H22").Value = "'***.*"
$ws.Range("C33").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("C27").Value = 1
$ws.Range("I33").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("F27").Value = 1
$ws.Range("I33").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("D31").Value = "'0"
$ws.Range("C33").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = "'***.*"
$ws.Range("C33").Copy()
$ws.Range("E31").PasteSpecial(-4122)

$excel.CutCopyMode = 0

